$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70 (ALC)
$ws.Range("H70").Value = 10480423
$ws.Range("I70").Value = 41917492
$ws.Range("J70").Value = 1399.75
$ws.Range("K70").Value = 125752476
$ws.Range("L70").Value = 4199.25
$ws.Range("M70").Value = -125752206
$ws.Range("N70").Value = -4739.25

# Row 73 (ALC)
$ws.Range("H73").Value = 10480423
$ws.Range("I73").Value = 41917492
$ws.Range("J73").Value = 1399.75
$ws.Range("K73").Value = 125752476
$ws.Range("L73").Value = 4199.25
$ws.Range("M73").Value = -125751540
$ws.Range("N73").Value = -6071.25

# Row 96 (ALC)
$ws.Range("H96").Value = 1033.3334
$ws.Range("J96").Value = 1300
$ws.Range("L96").Value = 3900
$ws.Range("N96").Value = -6646

# Row 111 (ALC)
$ws.Range("H111").Value = 900
$ws.Range("I111").Value = 900
$ws.Range("K111").Value = 2700
$ws.Range("M111").Value = 367

# Row 131 (ALC)
$ws.Range("H131").Value = 4636.12
$ws.Range("I131").Value = 1006.1667
$ws.Range("J131").Value = 5782.421
$ws.Range("K131").Value = 3018.5001
$ws.Range("L131").Value = 17347.263
$ws.Range("M131").Value = 2021.4999
$ws.Range("N131").Value = -27427.263

$ws = $wb.Worksheets.Item("ARM")
# Row 57 (ARM)
$ws.Range("H57").Value = 4989
$ws.Range("I57").Value = 4989
$ws.Range("K57").Value = 4989
$ws.Range("M57").Value = -4505

# Row 61 (ARM)
$ws.Range("H61").Value = 1590
$ws.Range("I61").Value = 1590
$ws.Range("K61").Value = 1590
$ws.Range("M61").Value = -1378

# Row 74 (ARM)
$ws.Range("H74").Value = 1105.4736
$ws.Range("I74").Value = 860.4
$ws.Range("J74").Value = 1377.7778
$ws.Range("K74").Value = 860.4
$ws.Range("L74").Value = 1377.7778
$ws.Range("M74").Value = 13.60000000000002
$ws.Range("N74").Value = -3125.7778

# Row 77 (ARM)
$ws.Range("H77").Value = 1105.4736
$ws.Range("I77").Value = 860.4
$ws.Range("J77").Value = 1377.7778
$ws.Range("K77").Value = 4302
$ws.Range("L77").Value = 6888.889
$ws.Range("M77").Value = 66
$ws.Range("N77").Value = -15624.889

# Row 132 (ARM)
$ws.Range("H132").Value = 1691.0625
$ws.Range("I132").Value = 1470.5333
$ws.Range("K132").Value = 4411.5999
$ws.Range("M132").Value = -1881.5999

# Row 135 (ARM)
$ws.Range("H135").Value = 28078.3
$ws.Range("J135").Value = 28078.3
$ws.Range("L135").Value = 28078.3
$ws.Range("N135").Value = -38218.3

# Row 136 (ARM)
$ws.Range("H136").Value = 1590
$ws.Range("I136").Value = 1590
$ws.Range("K136").Value = 4770
$ws.Range("M136").Value = -2220

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (BSM)
$ws.Range("H107").Value = 22611.885
$ws.Range("I107").Value = 26077.682
$ws.Range("J107").Value = 3550
$ws.Range("K107").Value = 26077.682
$ws.Range("L107").Value = 3550
$ws.Range("M107").Value = -24157.682
$ws.Range("N107").Value = -7390

$ws = $wb.Worksheets.Item("CRP")
# Row 28 (CRP)
$ws.Range("H28").Value = 188881
$ws.Range("J28").Value = 188881
$ws.Range("L28").Value = 188881
$ws.Range("N28").Value = -189371

# Row 31 (CRP)
$ws.Range("H31").Value = 34486544
$ws.Range("I31").Value = 55558790
$ws.Range("K31").Value = 55558790
$ws.Range("M31").Value = -55558495

# Row 34 (CRP)
$ws.Range("H34").Value = 34486544
$ws.Range("I34").Value = 55558790
$ws.Range("K34").Value = 55558790
$ws.Range("M34").Value = -55558588

# Row 76 (CRP)
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79 (CRP)
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 99 (CRP)
$ws.Range("H99").Value = 1997.92
$ws.Range("I99").Value = 2207.5789
$ws.Range("K99").Value = 2207.5789
$ws.Range("M99").Value = -709.5789

# Row 126 (CRP)
$ws.Range("H126").Value = 1997.92
$ws.Range("I126").Value = 2207.5789
$ws.Range("K126").Value = 6622.736699999999
$ws.Range("M126").Value = -4152.736699999999

# Row 131 (CRP)
$ws.Range("H131").Value = 22050.4
$ws.Range("I131").Value = 9296
$ws.Range("J131").Value = 25239
$ws.Range("K131").Value = 9296
$ws.Range("L131").Value = 25239
$ws.Range("N131").Value = -35319
$ws.Range("M131").Value = -4256

# Row 132 (CRP)
$ws.Range("H132").Value = 2704.318
$ws.Range("I132").Value = 2473.875
$ws.Range("J132").Value = 3318.8333
$ws.Range("K132").Value = 7421.625
$ws.Range("L132").Value = 9956.499899999999
$ws.Range("M132").Value = -4891.625
$ws.Range("N132").Value = -15016.4999

$ws = $wb.Worksheets.Item("CUL")
# Row 60 (CUL)
$ws.Range("H60").Value = 440.4
$ws.Range("I60").Value = 229.14285
$ws.Range("J60").Value = 933.3333
$ws.Range("K60").Value = 687.4285500000001
$ws.Range("L60").Value = 2799.9999
$ws.Range("M60").Value = -436.4285500000001
$ws.Range("N60").Value = -3301.9999

# Row 131 (CUL)
$ws.Range("H131").Value = 868.88776
$ws.Range("I131").Value = 474
$ws.Range("J131").Value = 890.1183
$ws.Range("K131").Value = 1422
$ws.Range("L131").Value = 2670.3549
$ws.Range("M131").Value = 3618
$ws.Range("N131").Value = -12750.3549

# Row 132 (CUL)
$ws.Range("H132").Value = 1886.5625
$ws.Range("I132").Value = 1101.2
$ws.Range("J132").Value = 2243.5454
$ws.Range("K132").Value = 9910.800000000001
$ws.Range("L132").Value = 20191.9086
$ws.Range("M132").Value = -7380.800000000001
$ws.Range("N132").Value = -25251.9086

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 3634.6086
$ws.Range("I122").Value = 4321.857
$ws.Range("K122").Value = 12965.571
$ws.Range("M122").Value = -10515.571

# Row 126 (GSM)
$ws.Range("H126").Value = 2417.95
$ws.Range("I126").Value = 2267.7144
$ws.Range("J126").Value = 2768.5
$ws.Range("K126").Value = 6803.1432
$ws.Range("L126").Value = 8305.5
$ws.Range("M126").Value = -4333.1432
$ws.Range("N126").Value = -13245.5

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Range("H61").Value = 2955
$ws.Range("I61").Value = 3100.0908
$ws.Range("J61").Value = 2635.8
$ws.Range("K61").Value = 3100.0908
$ws.Range("L61").Value = 2635.8
$ws.Range("M61").Value = -2898.0908
$ws.Range("N61").Value = -3039.8

# Row 93 (LTW)
$ws.Range("H93").Value = 1300.9333
$ws.Range("I93").Value = 1237.7273
$ws.Range("J93").Value = 1474.75
$ws.Range("K93").Value = 1237.7273
$ws.Range("L93").Value = 1474.75
$ws.Range("M93").Value = 10.27269999999999
$ws.Range("N93").Value = -3970.75

# Row 113 (LTW)
$ws.Range("H113").Value = 2955
$ws.Range("I113").Value = 3100.0908
$ws.Range("J113").Value = 2635.8
$ws.Range("K113").Value = 3100.0908
$ws.Range("L113").Value = 2635.8
$ws.Range("M113").Value = -930.0907999999999
$ws.Range("N113").Value = -6975.8

# Row 136 (LTW)
$ws.Range("H136").Value = 4635.515
$ws.Range("I136").Value = 5239.7036
$ws.Range("J136").Value = 1916.6666
$ws.Range("K136").Value = 15719.1108
$ws.Range("L136").Value = 5749.9998
$ws.Range("M136").Value = -13169.1108
$ws.Range("N136").Value = -10849.9998

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Range("H62").Value = 4277.778
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5748

# Row 65 (WVR)
$ws.Range("H65").Value = 4277.778
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -28740

# Row 113 (WVR)
$ws.Range("H113").Value = 760.2222
$ws.Range("I113").Value = 430.7
$ws.Range("J113").Value = 1701.7142
$ws.Range("K113").Value = 1292.1
$ws.Range("L113").Value = 5105.142599999999
$ws.Range("M113").Value = 877.9000000000001
$ws.Range("N113").Value = -9445.142599999999

# Row 132 (WVR)
$ws.Range("H132").Value = 2062.6191
$ws.Range("I132").Value = 2336.2354
$ws.Range("J132").Value = 899.75
$ws.Range("K132").Value = 7008.706200000001
$ws.Range("L132").Value = 2699.25
$ws.Range("M132").Value = -4478.706200000001
$ws.Range("N132").Value = -7759.25
